$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.732.25"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.946.04"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.39%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.945.63"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.15"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +10.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +9.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.56"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.432.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.687.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.946.52"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "438.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.666"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.44"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.93"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000104"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +24.25%  "
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("E33").Value = "  +6.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.20"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.13"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +12.23%  "
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.67"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  +7.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.42"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.280"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.95"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.706.28"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0342"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "358.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.83"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.75%  "
